$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 9999
$ws.Range("J10").Value = 9999
$ws.Range("L10").Value = 9999
$ws.Range("N10").Value = -10585
$ws.Range("H19").Value = 980.6429000000001
$ws.Range("I19").Value = 775.1818
$ws.Range("J19").Value = 1734
$ws.Range("K19").Value = 775.1818
$ws.Range("L19").Value = 1734
$ws.Range("M19").Value = -600.1818
$ws.Range("N19").Value = -2084
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = $null
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = $null
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = $null
$ws.Range("N77").Value = $null
$ws.Range("H112").Value = 999
$ws.Range("J112").Value = 999
$ws.Range("L112").Value = 2997
$ws.Range("N112").Value = -5213
$ws.Range("H125").Value = 1499
$ws.Range("J125").Value = 1499
$ws.Range("L125").Value = 13491
$ws.Range("N125").Value = -18411
$ws.Range("H129").Value = 1423.5
$ws.Range("I129").Value = 1347.5
$ws.Range("K129").Value = 4042.5
$ws.Range("M129").Value = 957.5
$ws.Range("H138").Value = 8311.450000000001
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("H141").Value = 5611.3335
$ws.Range("I141").Value = 5512.143
$ws.Range("K141").Value = 16536.429
$ws.Range("M141").Value = -11356.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1600
$ws.Range("I2").Value = 1700
$ws.Range("K2").Value = 1700
$ws.Range("M2").Value = -1587
$ws.Range("H74").Value = 1879
$ws.Range("I74").Value = 1866
$ws.Range("J74").Value = 1898.5
$ws.Range("K74").Value = 1866
$ws.Range("L74").Value = 1898.5
$ws.Range("M74").Value = -992
$ws.Range("N74").Value = -3646.5
$ws.Range("H77").Value = 1879
$ws.Range("I77").Value = 1866
$ws.Range("J77").Value = 1898.5
$ws.Range("K77").Value = 9330
$ws.Range("L77").Value = 9492.5
$ws.Range("M77").Value = -4962
$ws.Range("N77").Value = -18228.5
$ws.Range("H116").Value = 1600
$ws.Range("I116").Value = 1700
$ws.Range("K116").Value = 1700
$ws.Range("M116").Value = 594
$ws.Range("H132").Value = 2322.16
$ws.Range("I132").Value = 1812.238
$ws.Range("K132").Value = 5436.714
$ws.Range("M132").Value = -2906.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1600
$ws.Range("I3").Value = 1700
$ws.Range("K3").Value = 1700
$ws.Range("M3").Value = -1586
$ws.Range("H22").Value = 629
$ws.Range("I22").Value = 614.8
$ws.Range("K22").Value = 614.8
$ws.Range("M22").Value = -441.8
$ws.Range("H107").Value = 2050.7144
$ws.Range("I107").Value = 1495.4
$ws.Range("K107").Value = 1495.4
$ws.Range("M107").Value = 424.5999999999999
$ws.Range("H134").Value = 5400.3
$ws.Range("I134").Value = 5500.3335
$ws.Range("K134").Value = 16501.0005
$ws.Range("M134").Value = -13966.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5999.5
$ws.Range("I134").Value = 5999.5
$ws.Range("K134").Value = 17998.5
$ws.Range("M134").Value = -15463.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 300
$ws.Range("I51").Value = 300
$ws.Range("J51").Value = 300
$ws.Range("K51").Value = 900
$ws.Range("L51").Value = 900
$ws.Range("M51").Value = -440
$ws.Range("N51").Value = -1820
$ws.Range("H122").Value = 865.2222
$ws.Range("J122").Value = 941.1429000000001
$ws.Range("L122").Value = 8470.286100000001
$ws.Range("N122").Value = -13370.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 443.75
$ws.Range("I107").Value = 449
$ws.Range("J107").Value = 438.5
$ws.Range("K107").Value = 449
$ws.Range("L107").Value = 438.5
$ws.Range("M107").Value = 1471
$ws.Range("N107").Value = -4278.5
$ws.Range("H113").Value = 890
$ws.Range("J113").Value = 882.6667
$ws.Range("L113").Value = 882.6667
$ws.Range("N113").Value = -5222.6667
$ws.Range("H126").Value = 3774.6667
$ws.Range("I126").Value = 2529.8
$ws.Range("K126").Value = 7589.400000000001
$ws.Range("M126").Value = -5119.400000000001
$ws.Range("H132").Value = 3145
$ws.Range("I132").Value = 1624
$ws.Range("K132").Value = 4872
$ws.Range("M132").Value = -2342

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = $null
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = $null
$ws.Range("H132").Value = 4196.4165
$ws.Range("I132").Value = 3595.6667
$ws.Range("K132").Value = 10787.0001
$ws.Range("M132").Value = -8257.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 19999
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 19999
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 19999
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = -20229
$ws.Range("H70").Value = 46701.332
$ws.Range("I70").Value = 39999.5
$ws.Range("K70").Value = 39999.5
$ws.Range("M70").Value = -39684.5
$ws.Range("H73").Value = 46701.332
$ws.Range("I73").Value = 39999.5
$ws.Range("K73").Value = 39999.5
$ws.Range("M73").Value = -38907.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null
